$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state: table re-sorted by "realeffort" (column F) descending, using
# refreshed simulated realeffort values. Columns: B=index id, C=prolificid,
# D=name, E=gender, F=realeffort, G=race (constant "Hispanic", unchanged),
# H=re_rank (unchanged 1..12).

$rows = @(
    @{ Row=2;  B=7;  C='6024c18b094ac71dd93f4f5a'; D='Katherine'; E='female'; F=8.051697533201137 },
    @{ Row=3;  B=2;  C='60778ed0fde3e9c3a96f1d11'; D='Melissa';   E='female'; F=8.049345038247747 },
    @{ Row=4;  B=8;  C='5f0142aa1eb1e528e7abce50'; D='Valeria';   E='female'; F=7.172380869265427 },
    @{ Row=5;  B=3;  C='60ba8ba51a5e0a105396888a'; D='Alfredo';   E='male';   F=7.084258182079134 },
    @{ Row=6;  B=11; C='5f5ea8227fa75676f56f9276'; D='Carlos';    E='male';   F=6.251133082574972 },
    @{ Row=7;  B=0;  C='5eeaa065c7acf61c4322f6d9'; D='Yonifredy'; E='male';   F=6.228542514609791 },
    @{ Row=8;  B=6;  C='5e706891c396cc64388ef760'; D='Maria';     E='male';   F=5.257951082805501 },
    @{ Row=9;  B=4;  C='5dd671942b033b5ec8bc97b4'; D='Juan';      E='female'; F=3.374220793453376 },
    @{ Row=10; B=1;  C='60743a8fd12c5ffa72972fd5'; D='Mary';      E='female'; F=3.30059129783378 },
    @{ Row=11; B=10; C='5e0adc8f4cac6834756db412'; D='Josue';     E='male';   F=2.445237386649406 },
    @{ Row=12; B=12; C='5e58b3e415b8d40b5e1dabf1'; D='Cristian';  E='male';   F=1.109097411672669 },
    @{ Row=13; B=9;  C='5e35d91ea42bce592e996843'; D='Sergio';    E='male';   F=0.3238469706846104 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
}
